$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("B3").Value = 1.02
$ws.Range("B4").Value = 1.02
$ws.Range("B5").Value = 1.02
$ws.Range("B6").Value = 1.02
$ws.Range("B7").Value = 1.02
$ws.Range("B8").Value = 1.02
$ws.Range("B9").Value = 1.02
$ws.Range("B10").Value = 1.02
$ws.Range("B11").Value = 1.02
$ws.Range("B12").Value = 1.02
$ws.Range("B13").Value = 1.02
$ws.Range("B14").Value = 1.02
$ws.Range("B15").Value = 1.02
$ws.Range("B16").Value = 1.02
$ws.Range("B17").Value = 1.02
$ws.Range("B18").Value = 1.02
$ws.Range("B19").Value = 1.02
$ws.Range("B20").Value = 1.02
$ws.Range("B21").Value = 1.02
$ws.Range("B22").Value = 1.02
$ws.Range("B23").Value = 1.02
$ws.Range("B24").Value = 1.02
$ws.Range("B25").Value = 1.02
$ws.Range("C2").Value = 1.02318932764533
$ws.Range("C3").Value = 1.024018927336684
$ws.Range("C4").Value = 1.02455640810911
$ws.Range("C5").Value = 1.024782524823179
$ws.Range("C6").Value = 1.02482050013767
$ws.Range("C7").Value = 1.024559428865855
$ws.Range("C8").Value = 1.023469554145742
$ws.Range("C9").Value = 1.021554292175255
$ws.Range("C10").Value = 1.020281068812164
$ws.Range("C11").Value = 1.019730627613595
$ws.Range("C12").Value = 1.019526302031206
$ws.Range("C13").Value = 1.019570124526951
$ws.Range("C14").Value = 1.019713735278897
$ws.Range("C15").Value = 1.019802236247646
$ws.Range("C16").Value = 1.020317618313079
$ws.Range("C17").Value = 1.020641138729203
$ws.Range("C18").Value = 1.020829926828272
$ws.Range("C19").Value = 1.020894312908633
$ws.Range("C20").Value = 1.020606419338186
$ws.Range("C21").Value = 1.019671441813084
$ws.Range("C22").Value = 1.019084353414035
$ws.Range("C23").Value = 1.019395506633778
$ws.Range("C24").Value = 1.020622107276806
$ws.Range("C25").Value = 1.022048802742884
$ws.Range("D2").Value = 1.027559680384708
$ws.Range("D3").Value = 1.028138993805809
$ws.Range("D4").Value = 1.028514200357431
$ws.Range("D5").Value = 1.028672019665036
$ws.Range("D6").Value = 1.02869852300499
$ws.Range("D7").Value = 1.028516308824953
$ws.Range("D8").Value = 1.027755388094067
$ws.Range("D9").Value = 1.02641732433878
$ws.Range("D10").Value = 1.025527255712551
$ws.Range("D11").Value = 1.025142337111718
$ws.Range("D12").Value = 1.024999436216008
$ws.Range("D13").Value = 1.025030085519205
$ws.Range("D14").Value = 1.025130523338231
$ws.Range("D15").Value = 1.025192416388271
$ws.Range("D16").Value = 1.025552811910388
$ws.Range("D17").Value = 1.025779010185264
$ws.Range("D18").Value = 1.02591099467566
$ws.Range("D19").Value = 1.025956005868746
$ws.Range("D20").Value = 1.025754736390437
$ws.Range("D21").Value = 1.025100944818389
$ws.Range("D22").Value = 1.024690315010699
$ws.Range("D23").Value = 1.024907955801787
$ws.Range("D24").Value = 1.025765704528905
$ws.Range("D25").Value = 1.026762905847445
$ws.Range("E2").Value = 1.034127452226345
$ws.Range("E3").Value = 1.034894400584912
$ws.Range("E4").Value = 1.03539166785219
$ws.Range("E5").Value = 1.035600956686621
$ws.Range("E6").Value = 1.035636111087265
$ws.Range("E7").Value = 1.035394463448226
$ws.Range("E8").Value = 1.034386437781702
$ws.Range("E9").Value = 1.032617910138864
$ws.Range("E10").Value = 1.031444205946584
$ws.Range("E11").Value = 1.030937262510972
$ws.Range("E12").Value = 1.030749155050055
$ws.Range("E13").Value = 1.030789495934544
$ws.Range("E14").Value = 1.030921709523241
$ws.Range("E15").Value = 1.031003196428562
$ws.Range("E16").Value = 1.031477877202678
$ws.Range("E17").Value = 1.031775975444549
$ws.Range("E18").Value = 1.031949974310829
$ws.Range("E19").Value = 1.032009324268721
$ws.Range("E20").Value = 1.031743979567229
$ws.Range("E21").Value = 1.030882770549078
$ws.Range("E22").Value = 1.030342417041407
$ws.Range("E23").Value = 1.03062876161063
$ws.Range("E24").Value = 1.031758436746498
$ws.Range("E25").Value = 1.033074188211561
$ws.Range("F2").Value = 1.04591939632489
$ws.Range("F3").Value = 1.046839479773631
$ws.Range("F4").Value = 1.047436009194921
$ws.Range("F5").Value = 1.047687069281336
$ws.Range("F6").Value = 1.047729239686442
$ws.Range("F7").Value = 1.047439362777748
$ws.Range("F8").Value = 1.046230098615736
$ws.Range("F9").Value = 1.044108296253734
$ws.Range("F10").Value = 1.04269997770013
$ws.Range("F11").Value = 1.042091658795893
$ws.Range("F12").Value = 1.041865928226289
$ws.Range("F13").Value = 1.041914337934294
$ws.Range("F14").Value = 1.042072995211766
$ws.Range("F15").Value = 1.042170779223163
$ws.Range("F16").Value = 1.042740381430482
$ws.Range("F17").Value = 1.043098078761258
$ws.Range("F18").Value = 1.043306861419049
$ws.Range("F19").Value = 1.043378075251124
$ws.Range("F20").Value = 1.04305968632904
$ws.Range("F21").Value = 1.042026268339258
$ws.Range("F22").Value = 1.04137782687406
$ws.Range("F23").Value = 1.041721453191485
$ws.Range("F24").Value = 1.043077033772323
$ws.Range("F25").Value = 1.044655745733904
$ws.Range("I2").Value = 1.03048978317118
$ws.Range("I3").Value = 1.030593732677621
$ws.Range("I4").Value = 1.030659812627076
$ws.Range("I5").Value = 1.030687308920232
$ws.Range("I6").Value = 1.030691909012192
$ws.Range("I7").Value = 1.030660181149537
$ws.Range("I8").Value = 1.03052515780328
$ws.Range("I9").Value = 1.030278203588906
$ws.Range("I10").Value = 1.030107538693363
$ws.Range("I11").Value = 1.03003221733677
$ws.Range("I12").Value = 1.03000402660078
$ws.Range("I13").Value = 1.03001008323884
$ws.Range("I14").Value = 1.030029891423684
$ws.Range("I15").Value = 1.03004206769123
$ws.Range("I16").Value = 1.030112507613063
$ws.Range("I17").Value = 1.03015631238141
$ws.Range("I18").Value = 1.030181725683173
$ws.Range("I19").Value = 1.030190367657496
$ws.Range("I20").Value = 1.03015162673894
$ws.Range("I21").Value = 1.030024064282343
$ws.Range("I22").Value = 1.029942628500999
$ws.Range("I23").Value = 1.029985915727861
$ws.Range("I24").Value = 1.03015374440333
$ws.Range("I25").Value = 1.030343112691979
$ws.Range("J2").Value = 1.028370943506751
$ws.Range("J3").Value = 1.028839536690629
$ws.Range("J4").Value = 1.029142776754282
$ws.Range("J5").Value = 1.029270264549543
$ws.Range("J6").Value = 1.029291670607467
$ws.Range("J7").Value = 1.029144480231729
$ws.Range("J8").Value = 1.028529300020084
$ws.Range("J9").Value = 1.027445551054119
$ws.Range("J10").Value = 1.026723312962084
$ws.Range("J11").Value = 1.02641065364439
$ws.Range("J12").Value = 1.026294530522377
$ws.Range("J13").Value = 1.026319438741559
$ws.Range("J14").Value = 1.026401054609174
$ws.Range("J15").Value = 1.026451342529237
$ws.Range("J16").Value = 1.026744064810263
$ws.Range("J17").Value = 1.026927702617644
$ws.Range("J18").Value = 1.027034822550828
$ws.Range("J19").Value = 1.027071348846635
$ws.Range("J20").Value = 1.026907999274559
$ws.Range("J21").Value = 1.026377020415955
$ws.Range("J22").Value = 1.026043245707878
$ws.Range("J23").Value = 1.026220178673436
$ws.Range("J24").Value = 1.026916902345377
$ws.Range("J25").Value = 1.027725685578425
$ws.Range("K2").Value = 1.030379054823549
$ws.Range("K3").Value = 1.030766736607205
$ws.Range("K4").Value = 1.031017244454061
$ws.Range("K5").Value = 1.031122473324547
$ws.Range("K6").Value = 1.031140136715318
$ws.Range("K7").Value = 1.031018650861508
$ws.Range("K8").Value = 1.030510145028903
$ws.Range("K9").Value = 1.029611482218715
$ws.Range("K10").Value = 1.029010694264946
$ws.Range("K11").Value = 1.028750164113004
$ws.Range("K12").Value = 1.02865333505805
$ws.Range("K13").Value = 1.028674107757694
$ws.Range("K14").Value = 1.028742161339418
$ws.Range("K15").Value = 1.028784083941339
$ws.Range("K16").Value = 1.029027976791313
$ws.Range("K17").Value = 1.02918086216439
$ws.Range("K18").Value = 1.029270000382466
$ws.Range("K19").Value = 1.029300387856718
$ws.Range("K20").Value = 1.029164462852699
$ws.Range("K21").Value = 1.028722122806703
$ws.Range("K22").Value = 1.028443679481083
$ws.Range("K23").Value = 1.028591318083443
$ws.Range("K24").Value = 1.029171873110958
$ws.Range("K25").Value = 1.029844109648775
$ws.Range("L2").Value = 1.036927818908924
$ws.Range("L3").Value = 1.037504003075896
$ws.Range("L4").Value = 1.037877171933823
$ws.Range("L5").Value = 1.038034131880308
$ws.Range("L6").Value = 1.038060490781192
$ws.Range("L7").Value = 1.037879268928721
$ws.Range("L8").Value = 1.037122471951787
$ws.Range("L9").Value = 1.035791560344602
$ws.Range("L10").Value = 1.034906161197247
$ws.Range("L11").Value = 1.034523236508919
$ws.Range("L12").Value = 1.034381071650291
$ws.Range("L13").Value = 1.0344115632849
$ws.Range("L14").Value = 1.034511483672343
$ws.Range("L15").Value = 1.034573057295972
$ws.Range("L16").Value = 1.034931584442853
$ws.Range("L17").Value = 1.03515660296612
$ws.Range("L18").Value = 1.035287896651037
$ws.Range("L19").Value = 1.035332671855085
$ws.Range("L20").Value = 1.03513245602707
$ws.Range("L21").Value = 1.034482057652997
$ws.Range("L22").Value = 1.034073534535453
$ws.Range("L23").Value = 1.034290061192761
$ws.Range("L24").Value = 1.035143366853085
$ws.Range("L25").Value = 1.036135308659861
$ws.Range("M2").Value = 1.048686286123415
$ws.Range("M3").Value = 1.049417624984075
$ws.Range("M4").Value = 1.049891393302172
$ws.Range("M5").Value = 1.050090693851737
$ws.Range("M6").Value = 1.050124164768352
$ws.Range("M7").Value = 1.049894055863666
$ws.Range("M8").Value = 1.048933332032992
$ws.Range("M9").Value = 1.047244648631749
$ws.Range("M10").Value = 1.046121807606265
$ws.Range("M11").Value = 1.045636325496163
$ws.Range("M12").Value = 1.045456105028015
$ws.Range("M13").Value = 1.045494757960989
$ws.Range("M14").Value = 1.045621426178413
$ws.Range("M15").Value = 1.045699485169554
$ws.Range("M16").Value = 1.046154042675233
$ws.Range("M17").Value = 1.046439367348349
$ws.Range("M18").Value = 1.046605861319754
$ws.Range("M19").Value = 1.046662643098808
$ws.Range("M20").Value = 1.046408747570449
$ws.Range("M21").Value = 1.045584122517184
$ws.Range("M22").Value = 1.045066280684394
$ws.Range("M23").Value = 1.045340737895886
$ws.Range("M24").Value = 1.046422583117335
$ws.Range("M25").Value = 1.047680701557569
$ws.Range("N2").Value = 1.013506424914814
$ws.Range("N3").Value = 1.013663015106559
$ws.Range("N4").Value = 1.01376430557437
$ws.Range("N5").Value = 1.013806879518226
$ws.Range("N6").Value = 1.013814027354775
$ws.Range("N7").Value = 1.013764874483744
$ws.Range("N8").Value = 1.013559352035153
$ws.Range("N9").Value = 1.013196956390623
$ws.Range("N10").Value = 1.012955226330375
$ws.Range("N11").Value = 1.0128505288348
$ws.Range("N12").Value = 1.01281163596822
$ws.Range("N13").Value = 1.012819978773476
$ws.Range("N14").Value = 1.012847314008128
$ws.Range("N15").Value = 1.012864155691312
$ws.Range("N16").Value = 1.012962174234071
$ws.Range("N17").Value = 1.013023651836622
$ws.Range("N18").Value = 1.013059508066351
$ws.Range("N19").Value = 1.013071733656133
$ws.Range("N20").Value = 1.013017056138859
$ws.Range("N21").Value = 1.012839264557209
$ws.Range("N22").Value = 1.012727459596563
$ws.Range("N23").Value = 1.01278673129284
$ws.Range("N24").Value = 1.013020036458911
$ws.Range("N25").Value = 1.013290669807185
